$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10-13 (logo-web-social, servicios-diseno, black-friday,
# tarjetas-navidad) - the remaining "calendarios-regalos" row shifts up
# to become the new row 10, keeping its original display_order value.
$ws.Range("A10:G13").EntireRow.Delete()

# Update the active selection to match the post-edit state.
$ws.Range("A10:XFD10").Select()

# Match the window position recorded for the save.
$excel.ActiveWindow.Top = 110
